$wb = $excel.ActiveWorkbook

# --- Original "Sheet1": add the attendance header row and leave the
#     selection on C1 (matches xl/worksheets/sheet3.xml in the target) ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("A1").Value = "STUDENT"
$sheet1.Range("B1").Value = "DATE"
$sheet1.Range("C1").Value = "TIME"
$sheet1.Range("C1").Select() | Out-Null

# --- New sheet for 2024-01-07, inserted right before "Sheet1" (which is
#     the active sheet at this point) and populated with the same header
#     row; selection left on B5 (matches xl/worksheets/sheet2.xml) ---
$jan07 = $wb.Worksheets.Add()
$jan07.Name = "2024-01-07"
$jan07.Range("A1").Value = "STUDENT"
$jan07.Range("B1").Value = "DATE"
$jan07.Range("C1").Value = "TIME"
$jan07.Range("B5").Select() | Out-Null

# --- New sheet for 2024-01-08, inserted right before "2024-01-07" (now
#     the active sheet); left empty, becomes the active/selected tab with
#     the cursor on C1 (matches xl/worksheets/sheet1.xml) ---
$jan08 = $wb.Worksheets.Add()
$jan08.Name = "2024-01-08"
$jan08.Range("C1").Select() | Out-Null
